$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $range = $d.Content
    $range.Find.Execute($old, $true, $true, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2024-08-12 Monday" "2024-08-13 Tuesday"

Replace-Text "850×8=6800" "832×5=4160"
Replace-Text "145×6=870" "556×4=2224"
Replace-Text "737×4=2948" "231×9=2079"
Replace-Text "949×2=1898" "462×3=1386"
Replace-Text "964×3=2892" "552×4=2208"
Replace-Text "927×6=5562" "548×6=3288"
Replace-Text "358×7=2506" "943×7=6601"
Replace-Text "537×2=1074" "817×6=4902"
Replace-Text "154×4=616" "577×8=4616"
Replace-Text "281×7=1967" "182×3=546"
Replace-Text "975×2=1950" "107×9=963"
Replace-Text "997×9=8973" "899×2=1798"
Replace-Text "363×6=2178" "301×3=903"
Replace-Text "631×3=1893" "453×3=1359"
Replace-Text "578×9=5202" "834×4=3336"
Replace-Text "626×4=2504" "534×7=3738"
Replace-Text "920×2=1840" "249×3=747"
Replace-Text "152×3=456" "158×5=790"
Replace-Text "976×8=7808" "765×9=6885"
Replace-Text "864×4=3456" "870×6=5220"
Replace-Text "285×9=2565" "505×9=4545"
Replace-Text "376×5=1880" "993×5=4965"
Replace-Text "270×6=1620" "927×9=8343"
Replace-Text "296×4=1184" "368×8=2944"
Replace-Text "878×8=7024" "389×8=3112"
